$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.338145017623901
$ws.Range("B1").Value = 1.993512153625488
$ws.Range("C1").Value = 3.806264877319336
$ws.Range("D1").Value = 0.9684718251228333
$ws.Range("E1").Value = 0.7653455138206482
